$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated power-flow line results (pl_mw) for the 380 kV case - rows 2-25, columns B:L
$data = New-Object "object[,]" 24,11
$data[0,0] = 2.686394174516352
$data[0,1] = 0.3325364391967014
$data[0,2] = 0.01496881308703024
$data[0,3] = 0
$data[0,4] = 3.963262293284828
$data[0,5] = 0.002605820754676929
$data[0,6] = 0
$data[0,7] = 2.331366034889314
$data[0,8] = 0.1380561296583309
$data[0,9] = 0
$data[0,10] = 0.4776836227994323
$data[1,0] = 2.593448375293065
$data[1,1] = 0.3040863828185536
$data[1,2] = 0.01432943959175859
$data[1,3] = 0
$data[1,4] = 3.932235353567819
$data[1,5] = 0.002611799434967612
$data[1,6] = 0
$data[1,7] = 2.320481800466922
$data[1,8] = 0.1386378639724875
$data[1,9] = 0
$data[1,10] = 0.4724143378576855
$data[2,0] = 2.538031024753877
$data[2,1] = 0.2868107735639853
$data[2,2] = 0.01393273769992476
$data[2,3] = 0
$data[2,4] = 3.915172326769138
$data[2,5] = 0.00261566246422732
$data[2,6] = 0
$data[2,7] = 2.314890436770369
$data[2,8] = 0.1390222083820163
$data[2,9] = 0
$data[2,10] = 0.4694227632581232
$data[3,0] = 2.515862969610112
$data[3,1] = 0.2798187316864187
$data[3,2] = 0.01376998595525336
$data[3,3] = 0
$data[3,4] = 3.908716594195283
$data[3,5] = 0.002617285154698622
$data[3,6] = 0
$data[3,7] = 2.312885097919164
$data[3,8] = 0.1391856628037953
$data[3,9] = 0
$data[3,10] = 0.4682649750863135
$data[4,0] = 2.512207033892992
$data[4,1] = 0.2786605852759294
$data[4,2] = 0.013742893549054
$data[4,3] = 0
$data[4,4] = 3.907674613363795
$data[4,5] = 0.002617557533891728
$data[4,6] = 0
$data[4,7] = 2.312568574746962
$data[4,8] = 0.1392132169202522
$data[4,9] = 0
$data[4,10] = 0.4680764278940899
$data[5,0] = 2.537730378919207
$data[5,1] = 0.2867162832164354
$data[5,2] = 0.0139305472685507
$data[5,3] = 0
$data[5,4] = 3.915083250749632
$data[5,5] = 0.002615684152060593
$data[5,6] = 0
$data[5,7] = 2.314862287695348
$data[5,8] = 0.1390243851244382
$data[5,9] = 0
$data[5,10] = 0.4694069007047119
$data[6,0] = 2.654003361982745
$data[6,1] = 0.3226864731625483
$data[6,2] = 0.01474917946050169
$data[6,3] = 0
$data[6,4] = 3.952150312157841
$data[6,5] = 0.002607842437203354
$data[6,6] = 0
$data[6,7] = 2.32738579464413
$data[6,8] = 0.1382510764265703
$data[6,9] = 0
$data[6,10] = 0.475816171178792
$data[7,0] = 2.895155656851728
$data[7,1] = 0.3947844985885354
$data[7,2] = 0.01632434143657235
$data[7,3] = 0
$data[7,4] = 4.040718067152199
$data[7,5] = 0.002593981249548044
$data[7,6] = 0
$data[7,7] = 2.360668108777787
$data[7,8] = 0.1369500070202943
$data[7,9] = 0
$data[7,10] = 0.4903205719077448
$data[8,0] = 3.080417198481996
$data[8,1] = 0.4487546330212808
$data[8,2] = 0.01746681298327601
$data[8,3] = 0
$data[8,4] = 4.115631242476553
$data[8,5] = 0.002584710931890057
$data[8,6] = 0
$data[8,7] = 2.390530414718583
$data[8,8] = 0.1361252873603789
$data[8,9] = 0
$data[8,10] = 0.5021613588273368
$data[9,0] = 3.166472353199595
$data[9,1] = 0.4735351063784492
$data[9,2] = 0.01798412248440684
$data[9,3] = 0
$data[9,4] = 4.151884406727874
$data[9,5] = 0.00258068964355903
$data[9,6] = 0
$data[9,7] = 2.40531045252871
$data[9,8] = 0.1357785592478216
$data[9,9] = 0
$data[9,10] = 0.5078063618702657
$data[10,0] = 3.199316115465422
$data[10,1] = 0.4829525747431944
$data[10,2] = 0.01817973088228442
$data[10,3] = 0
$data[10,4] = 4.16592801372181
$data[10,5] = 0.002579194866547108
$data[10,6] = 0
$data[10,7] = 2.411080749363549
$data[10,8] = 0.1356513508118731
$data[10,9] = 0
$data[10,10] = 0.5099812220407642
$data[11,0] = 3.192231201003381
$data[11,1] = 0.4809228471905271
$data[11,2] = 0.01813761483994725
$data[11,3] = 0
$data[11,4] = 4.162889406137651
$data[11,5] = 0.002579515550889526
$data[11,6] = 0
$data[11,7] = 2.409830275827673
$data[11,8] = 0.1356785655177415
$data[11,9] = 0
$data[11,10] = 0.5095111709553919
$data[12,0] = 3.169169283205065
$data[12,1] = 0.4743092089824472
$data[12,2] = 0.01800022067763507
$data[12,3] = 0
$data[12,4] = 4.153033448292547
$data[12,5] = 0.002580566106934668
$data[12,6] = 0
$data[12,7] = 2.405781694125011
$data[12,8] = 0.1357680117624049
$data[12,9] = 0
$data[12,10] = 0.5079845427429177
$data[13,0] = 3.155076625762149
$data[13,1] = 0.4702625662160926
$data[13,2] = 0.01791602741500498
$data[13,3] = 0
$data[13,4] = 4.147037530946534
$data[13,5] = 0.002581213244271088
$data[13,6] = 0
$data[13,7] = 2.403324449816651
$data[13,8] = 0.1358233328130485
$data[13,9] = 0
$data[13,10] = 0.5070542872518331
$data[14,0] = 3.074829148319111
$data[14,1] = 0.4471398407789025
$data[14,2] = 0.01743296233155078
$data[14,3] = 0
$data[14,4] = 4.113305995757372
$data[14,5] = 0.002584977659219279
$data[14,6] = 0
$data[14,7] = 2.389588689585821
$data[14,8] = 0.1361485191779277
$data[14,9] = 0
$data[14,10] = 0.5017976508122359
$data[15,0] = 3.026056021914201
$data[15,1] = 0.433013980112662
$data[15,2] = 0.01713604410655378
$data[15,3] = 0
$data[15,4] = 4.093171638169935
$data[15,5] = 0.002587337045996091
$data[15,6] = 0
$data[15,7] = 2.38146951242615
$data[15,8] = 0.136355295010647
$data[15,9] = 0
$data[15,10] = 0.4986391258487117
$data[16,0] = 2.998170384608727
$data[16,1] = 0.4249106891483621
$data[16,2] = 0.01696503442909858
$data[16,3] = 0
$data[16,4] = 4.081795432622471
$data[16,5] = 0.002588712544760947
$data[16,6] = 0
$data[16,7] = 2.376912014280066
$data[16,8] = 0.1364769039861464
$data[16,9] = 0
$data[16,10] = 0.4968467610252105
$data[17,0] = 2.988757511616313
$data[17,1] = 0.4221707313573688
$data[17,2] = 0.01690709212255115
$data[17,3] = 0
$data[17,4] = 4.07797869719343
$data[17,5] = 0.002589181436996602
$data[17,6] = 0
$data[17,7] = 2.37538818946723
$data[17,8] = 0.1365185384657401
$data[17,9] = 0
$data[17,10] = 0.4962440756152233
$data[18,0] = 3.031230676500513
$data[18,1] = 0.4345154688671187
$data[18,2] = 0.01716767494335159
$data[18,3] = 0
$data[18,4] = 4.095293786691769
$data[18,5] = 0.002587083977761739
$data[18,6] = 0
$data[18,7] = 2.382322164977765
$data[18,8] = 0.1363330063115082
$data[18,9] = 0
$data[18,10] = 0.4989728373650877
$data[19,0] = 3.17593615861341
$data[19,1] = 0.4762508763096776
$data[19,2] = 0.01804058393700103
$data[19,3] = 0
$data[19,4] = 4.155919803520135
$data[19,5] = 0.002580256774683516
$data[19,6] = 0
$data[19,7] = 2.406966142166794
$data[19,8] = 0.1357416282288497
$data[19,9] = 0
$data[19,10] = 0.5084319399838506
$data[20,0] = 3.272005282623638
$data[20,1] = 0.503723791121331
$data[20,2] = 0.01860944490441696
$data[20,3] = 0
$data[20,4] = 4.197381564022606
$data[20,5] = 0.00257595791980133
$data[20,6] = 0
$data[20,7] = 2.42408392246341
$data[20,8] = 0.1353789681298654
$data[20,9] = 0
$data[20,10] = 0.5148309787405339
$data[21,0] = 3.220594256123036
$data[21,1] = 0.4890427776240927
$data[21,2] = 0.01830596259446082
$data[21,3] = 0
$data[21,4] = 4.17508349549766
$data[21,5] = 0.002578237427757677
$data[21,6] = 0
$data[21,7] = 2.414854786823057
$data[21,8] = 0.1355703452305068
$data[21,9] = 0
$data[21,10] = 0.5113958247120678
$data[22,0] = 3.028890733416119
$data[22,1] = 0.4338365902394798
$data[22,2] = 0.01715337560342789
$data[22,3] = 0
$data[22,4] = 4.094333742750308
$data[22,5] = 0.002587198330697591
$data[22,6] = 0
$data[22,7] = 2.381936337235018
$data[22,8] = 0.1363430745269572
$data[22,9] = 0
$data[22,10] = 0.4988218931529929
$data[23,0] = 2.828503107112226
$data[23,1] = 0.3751080563264679
$data[23,2] = 0.01590107882192626
$data[23,3] = 0
$data[23,4] = 4.015041509852608
$data[23,5] = 0.002597569860662342
$data[23,6] = 0
$data[23,7] = 2.350721005909691
$data[23,8] = 0.1372789358172337
$data[23,9] = 0
$data[23,10] = 0.4861890806642464

$ws.Range("B2:L25").Value = $data
